$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 (Exhibition) - F column ("want to go" count) updates ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(2,6).Value = 6821
$ws1.Cells.Item(4,6).Value = 148
$ws1.Cells.Item(5,6).Value = 16
$ws1.Cells.Item(6,6).Value = 755
$ws1.Cells.Item(7,6).Value = 755
$ws1.Cells.Item(8,6).Value = 24
$ws1.Cells.Item(10,6).Value = 29
$ws1.Cells.Item(11,6).Value = 1119
$ws1.Cells.Item(13,6).Value = 16
$ws1.Cells.Item(14,6).Value = 707
$ws1.Cells.Item(15,6).Value = 1022
$ws1.Cells.Item(16,6).Value = 1371
$ws1.Cells.Item(17,6).Value = 50
$ws1.Cells.Item(19,6).Value = 1548
$ws1.Cells.Item(20,6).Value = 8
$ws1.Cells.Item(21,6).Value = 586
$ws1.Cells.Item(26,6).Value = 1499
$ws1.Cells.Item(27,6).Value = 732
$ws1.Cells.Item(28,6).Value = 562
$ws1.Cells.Item(29,6).Value = 474
$ws1.Cells.Item(31,6).Value = 97
$ws1.Cells.Item(32,6).Value = 1012
$ws1.Cells.Item(33,6).Value = 1137
$ws1.Cells.Item(34,6).Value = 280
$ws1.Cells.Item(35,6).Value = 2377
$ws1.Cells.Item(37,6).Value = 1293
$ws1.Cells.Item(40,6).Value = 3901

# ---- Sheet: 演出 (Performance) ----
# Row 2 (2024-03-09 event) is removed; all subsequent rows (3..33) shift up by
# one position in columns B:I (column A serial numbers stay put), then the now-
# empty last row (33) is deleted, shrinking the sheet from A1:I33 to A1:I32.
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("B3:I33").Copy($ws2.Range("B2:I32"))
$ws2.Rows.Item(33).Delete()

# A handful of rows also got an updated "want to go" count as part of the refresh
$ws2.Cells.Item(6,6).Value = 9
$ws2.Cells.Item(7,6).Value = 170
$ws2.Cells.Item(13,6).Value = 395
$ws2.Cells.Item(14,6).Value = 340
$ws2.Cells.Item(18,6).Value = 24
$ws2.Cells.Item(19,6).Value = 40
$ws2.Cells.Item(21,6).Value = 245
$ws2.Cells.Item(23,6).Value = 114
$ws2.Cells.Item(25,6).Value = 230

# ---- Sheet: 本地生活 (Local life) - F column updates ----
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(5,6).Value = 1648
$ws3.Cells.Item(8,6).Value = 977

# ---- Sheet: 全部类型 (All types) - F column updates ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(4,6).Value = 1648
$ws4.Cells.Item(7,6).Value = 977
$ws4.Cells.Item(8,6).Value = 6821
$ws4.Cells.Item(11,6).Value = 148
$ws4.Cells.Item(12,6).Value = 16
$ws4.Cells.Item(13,6).Value = 755
$ws4.Cells.Item(14,6).Value = 755
$ws4.Cells.Item(16,6).Value = 29
$ws4.Cells.Item(17,6).Value = 1119
$ws4.Cells.Item(19,6).Value = 707
$ws4.Cells.Item(20,6).Value = 170
$ws4.Cells.Item(21,6).Value = 170
$ws4.Cells.Item(23,6).Value = 1022
$ws4.Cells.Item(24,6).Value = 1371
$ws4.Cells.Item(25,6).Value = 50
$ws4.Cells.Item(27,6).Value = 1548
$ws4.Cells.Item(28,6).Value = 8
$ws4.Cells.Item(29,6).Value = 586
$ws4.Cells.Item(31,6).Value = 340
$ws4.Cells.Item(34,6).Value = 1499
$ws4.Cells.Item(35,6).Value = 732
$ws4.Cells.Item(36,6).Value = 562
$ws4.Cells.Item(37,6).Value = 474
$ws4.Cells.Item(39,6).Value = 97
$ws4.Cells.Item(40,6).Value = 40
$ws4.Cells.Item(42,6).Value = 1012
$ws4.Cells.Item(43,6).Value = 1137
$ws4.Cells.Item(44,6).Value = 280
$ws4.Cells.Item(45,6).Value = 2377
$ws4.Cells.Item(46,6).Value = 230
$ws4.Cells.Item(49,6).Value = 1293
$ws4.Cells.Item(51,6).Value = 3901
